$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Genotypes" (col B) and "exp" (col D) swap places.
$headerGenotype = $ws.Cells.Item(1, 2).Value2
$headerExp = $ws.Cells.Item(1, 4).Value2
$ws.Cells.Item(1, 2).Value = $headerExp
$ws.Cells.Item(1, 4).Value = $headerGenotype

# The sheet is being restructured: the "Genotypes" values that used to live in
# column B are moved over to column D (which used to hold a redundant copy of
# the experiment/treatment stage). Column B is repurposed to hold a
# generation/stage label - "P" for the first block of rows, "M" for the
# middle block, and "H" for the last block.
for ($r = 2; $r -le 93; $r++) {
    $genotype = $ws.Cells.Item($r, 2).Value2

    if ($r -le 47) {
        $label = "P"
    } elseif ($r -le 77) {
        $label = "M"
    } else {
        $label = "H"
    }

    $ws.Cells.Item($r, 4).Value = $genotype
    $ws.Cells.Item($r, 2).Value = $label
}

# Column width tweaks that came along with the restructuring: column D (now
# holding the short genotype codes) is narrowed from 8 to 7 characters.
$ws.Columns.Item(4).ColumnWidth = 6.142857142857143

# Scroll/selection state, matching the saved view in the edited workbook.
$ws.Range("C76").Select()
$excel.ActiveWindow.ScrollRow = 59
$excel.ActiveWindow.ScrollColumn = 1
